# Update profit-calculation figures on Sheets (ALC/ARM/BSM/CRP/GSM/LTW/WVR)
# per scheduled-runner refresh of Gilgamesh_Profits source data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 799
$ws.Range("I33").Value = 799
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 799
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -570
$ws.Range("N33").ClearContents()
# Row 40
$ws.Range("H40").Value = 7686.25
$ws.Range("I40").Value = 7098.4
$ws.Range("J40").Value = 8666
$ws.Range("K40").Value = 7098.4
$ws.Range("L40").Value = 8666
$ws.Range("M40").Value = -6923.4
$ws.Range("N40").Value = -9016
# Row 98
$ws.Range("H98").Value = 1912.8667
$ws.Range("I98").Value = 1947.262
$ws.Range("K98").Value = 1947.262
$ws.Range("M98").Value = -449.2619999999999
# Row 101
$ws.Range("H101").Value = 400.5
$ws.Range("I101").Value = 352
$ws.Range("J101").Value = 497.5
$ws.Range("K101").Value = 1056
$ws.Range("L101").Value = 1492.5
$ws.Range("M101").Value = 566
$ws.Range("N101").Value = -4736.5
# Row 122
$ws.Range("H122").Value = 1912.8667
$ws.Range("I122").Value = 1947.262
$ws.Range("K122").Value = 5841.786
$ws.Range("M122").Value = -3391.786
# Row 132
$ws.Range("H132").Value = 4002.9836
$ws.Range("I132").Value = 4202.8335
$ws.Range("J132").Value = 2461.2856
$ws.Range("K132").Value = 12608.5005
$ws.Range("L132").Value = 7383.8568
$ws.Range("M132").Value = -10078.5005
$ws.Range("N132").Value = -12443.8568
# Row 138
$ws.Range("H138").Value = 3369.077
$ws.Range("I138").Value = 3067.56
$ws.Range("K138").Value = 9202.68
$ws.Range("M138").Value = -4062.68

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 869.36365
$ws.Range("I2").Value = 784.44446
$ws.Range("J2").Value = 1251.5
$ws.Range("K2").Value = 784.44446
$ws.Range("L2").Value = 1251.5
$ws.Range("M2").Value = -671.44446
$ws.Range("N2").Value = -1477.5
# Row 32
$ws.Range("H32").Value = 2902.9656
$ws.Range("I32").Value = 2902.9656
$ws.Range("K32").Value = 2902.9656
$ws.Range("M32").Value = -2615.9656
# Row 61
$ws.Range("H61").Value = 2228.96
$ws.Range("I61").Value = 2035.55
$ws.Range("J61").Value = 3002.6
$ws.Range("K61").Value = 2035.55
$ws.Range("L61").Value = 3002.6
$ws.Range("M61").Value = -1823.55
$ws.Range("N61").Value = -3426.6
# Row 74
$ws.Range("H74").Value = 310198.16
$ws.Range("I74").Value = 371237.2
$ws.Range("J74").Value = 5003
$ws.Range("K74").Value = 371237.2
$ws.Range("L74").Value = 5003
$ws.Range("M74").Value = -370363.2
$ws.Range("N74").Value = -6751
# Row 77
$ws.Range("H77").Value = 310198.16
$ws.Range("I77").Value = 371237.2
$ws.Range("J77").Value = 5003
$ws.Range("K77").Value = 1856186
$ws.Range("L77").Value = 25015
$ws.Range("M77").Value = -1851818
$ws.Range("N77").Value = -33751
# Row 97
$ws.Range("H97").Value = 1366
$ws.Range("I97").Value = 1142.7693
$ws.Range("K97").Value = 1142.7693
$ws.Range("M97").Value = -646.7692999999999
# Row 116
$ws.Range("H116").Value = 869.36365
$ws.Range("I116").Value = 784.44446
$ws.Range("J116").Value = 1251.5
$ws.Range("K116").Value = 784.44446
$ws.Range("L116").Value = 1251.5
$ws.Range("M116").Value = 1509.55554
$ws.Range("N116").Value = -5839.5
# Row 132
$ws.Range("H132").Value = 5850449.5
$ws.Range("I132").Value = 2111.0715
$ws.Range("J132").Value = 22225796
$ws.Range("K132").Value = 6333.2145
$ws.Range("L132").Value = 66677388
$ws.Range("M132").Value = -3803.2145
$ws.Range("N132").Value = -66682448
# Row 136
$ws.Range("H136").Value = 2228.96
$ws.Range("I136").Value = 2035.55
$ws.Range("J136").Value = 3002.6
$ws.Range("K136").Value = 6106.65
$ws.Range("L136").Value = 9007.799999999999
$ws.Range("M136").Value = -3556.65
$ws.Range("N136").Value = -14107.8

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 869.36365
$ws.Range("I3").Value = 784.44446
$ws.Range("J3").Value = 1251.5
$ws.Range("K3").Value = 784.44446
$ws.Range("L3").Value = 1251.5
$ws.Range("M3").Value = -670.44446
$ws.Range("N3").Value = -1479.5
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 134
$ws.Range("H134").Value = 2654.44
$ws.Range("I134").Value = 2334.7046
$ws.Range("K134").Value = 7004.1138
$ws.Range("M134").Value = -4469.1138

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1842.5
$ws.Range("I16").Value = 1946
$ws.Range("J16").Value = 1463
$ws.Range("K16").Value = 1946
$ws.Range("L16").Value = 1463
$ws.Range("M16").Value = -1659
$ws.Range("N16").Value = -2037
# Row 31
$ws.Range("H31").Value = 5262.418
$ws.Range("I31").Value = 4315.119
$ws.Range("K31").Value = 4315.119
$ws.Range("M31").Value = -4020.119
# Row 34
$ws.Range("H34").Value = 5262.418
$ws.Range("I34").Value = 4315.119
$ws.Range("K34").Value = 4315.119
$ws.Range("M34").Value = -4113.119
# Row 94
$ws.Range("H94").Value = 1700.8823
$ws.Range("I94").Value = 1580.7142
$ws.Range("J94").Value = 1785
$ws.Range("K94").Value = 1580.7142
$ws.Range("L94").Value = 1785
$ws.Range("M94").Value = -1129.7142
$ws.Range("N94").Value = -2687
# Row 113
$ws.Range("H113").Value = 1842.5
$ws.Range("I113").Value = 1946
$ws.Range("J113").Value = 1463
$ws.Range("K113").Value = 1946
$ws.Range("L113").Value = 1463
$ws.Range("M113").Value = 224
$ws.Range("N113").Value = -5803
# Row 132
$ws.Range("H132").Value = 23150390
$ws.Range("I132").Value = 16130826
$ws.Range("J132").Value = 66671690
$ws.Range("K132").Value = 48392478
$ws.Range("L132").Value = 200015070
$ws.Range("M132").Value = -48389948
$ws.Range("N132").Value = -200020130

$ws = $wb.Worksheets.Item("GSM")
# Row 24
$ws.Range("H24").Value = 10542.071
$ws.Range("J24").Value = 10466.333
$ws.Range("L24").Value = 10466.333
$ws.Range("N24").Value = -10812.333
# Row 132
$ws.Range("H132").Value = 1831.2222
$ws.Range("I132").Value = 1385.8334
$ws.Range("J132").Value = 4058.1667
$ws.Range("K132").Value = 4157.5002
$ws.Range("L132").Value = 12174.5001
$ws.Range("M132").Value = -1627.5002
$ws.Range("N132").Value = -17234.5001

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2971.6667
$ws.Range("I132").Value = 2377.2222
$ws.Range("K132").Value = 7131.6666
$ws.Range("M132").Value = -4601.6666
# Row 136
$ws.Range("H136").Value = 3076.0356
$ws.Range("I136").Value = 2975.2173
$ws.Range("K136").Value = 8925.651899999999
$ws.Range("M136").Value = -6375.651899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 868.57574
$ws.Range("I113").Value = 675.7083
$ws.Range("J113").Value = 1382.8889
$ws.Range("K113").Value = 2027.1249
$ws.Range("L113").Value = 4148.6667
$ws.Range("M113").Value = 142.8751
$ws.Range("N113").Value = -8488.6667
# Row 132
$ws.Range("H132").Value = 9012415
$ws.Range("I132").Value = 12348852
$ws.Range("J132").Value = 4034.3
$ws.Range("K132").Value = 37046556
$ws.Range("L132").Value = 12102.9
$ws.Range("M132").Value = -37044026
$ws.Range("N132").Value = -17162.9
